# Apply the "Add cantrals by cantons" edit:
#  - remove the old "units" row (row 2) and the stray no-index data row (old row 16)
#  - replace the old two-row header (row1 + row2) with a single header row
#    containing idx / idx2 / Name / Date Start / Date End / (m3/s) / (MW1) / (MW2) /
#    (GWh) Winter / (GWh) Summer / (GWh) Year
#  - tidy up the window/selection state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the stray data row (old row 16, which has no idx/name, just 3 numbers)
# first -- removing it before the units row keeps the row numbers stable while we
# reason about them.
$ws.Rows.Item(16).Delete()

# Remove the old "units" row (was row 2, directly under the header labels).
$ws.Rows.Item(2).Delete()

# Rebuild row 1 as the single header row for the new table layout.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 keep the default style; F1:K1 use the "text, Arial 9, no explicit number
# format" style (same font as the rest of the table headers used to have).
$ws.Range("A1:E1").Style = "Normal"

$headerStyle = $wb.Styles.Add("HeaderStyle")
$headerStyle.Font.Size = 9
$headerStyle.Font.Name = "Arial"
$ws.Range("F1:K1").Style = "HeaderStyle"
$headerStyle.Delete()

# Restore the view state recorded in the saved workbook.
$ws.Range("A15:K15").Select()

$wb.Windows.Item(1).WindowState = -4143
